$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4
    3  = 1
    4  = 2
    5  = 2
    6  = 4
    7  = 5
    8  = 5
    9  = 5
    10 = 3
    11 = 2
    12 = 4
    13 = 1
    14 = 3
    15 = 6
    16 = 3
    17 = 3
    18 = 5
    19 = 7
    20 = 2
    21 = 2
    22 = 3
    23 = 2
    24 = 2
    25 = 3
    26 = 1
    27 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
